# Apply the TestCases.xlsx content corrections (typo fixes) and update the
# saved sheet selection/scroll position, matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix small typos / wording issues in several test-case cells ---

# E4: "correclty" -> "correctly"
$ws.Range("E4").Value = "Access the application and validate that search field and ""Find by name"" button is displaying correctly"

# C6: "hompage" -> "homepage"
$ws.Range("C6").Value = "Validate that ""Add a new computer"" button is displaying on the homepage of application."

# C10: "Alphabatical" -> "Alphabetical"
$ws.Range("C10").Value = "Validate that member can search the computer name by entering the computer name (Alphabetical, String, special char etc)"

# C32: "infomration" -> "information"
$ws.Range("C32").Value = "Validate that application is displaying the error message for mandatory information if member tries to create a computer without entering any information"

# E51: "Validae" -> "Validate"
$ws.Range("E51").Value = "Validate that application populates all the existing information"

# H53: "retriving" -> "retrieving"
$ws.Range("H53").Value = "By including this scenario in the regression testing we can ensure that:" + [char]10 + "1. Member can update the computer information" + [char]10 + "2. All the existing information is being pre-populated. This will prove how we are caching the data/retrieving the data from DB while prepopulating" + [char]10 + "3. In other way this case will also prove member can add computer record with all the information"

# --- Update the sheet's saved view (scroll position / active selection) ---
$window = $excel.ActiveWindow
$window.ScrollRow = $ws.Range("E40").Row
$window.ScrollColumn = $ws.Range("E40").Column
$ws.Range("H50").Select()
